# Adding generated reports from pipeline run
# Insert a new "Tool Category" column after the "Tool" column on the
# Summary sheet, shifting the existing severity columns (INFO..UNKNOWN)
# one position to the right, and populate the new column's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Insert a new column at B; this shifts B:G -> C:H and preserves
# formatting/styles of the shifted cells (including the header style).
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Tool Category"

# Values for the newly inserted column, one per tool row.
$ws.Range("B2").Value = "SAST"
$ws.Range("B3").Value = "SCA"
$ws.Range("B4").Value = "IaC Scan"
$ws.Range("B5").Value = "Secret Scan"
